# Sync attendance_reports: reorder "Recorded By" entries in column G
# so that the first comma-separated name/email is moved to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -and $val.ToString().Contains(",")) {
        $parts = $val.ToString().Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $first = $parts[0]
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + @($first)
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value2 = $newVal
    }
}
